$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)

# Try setting error bars using a Range union like C3,E3,G3,I3,K3,M3,O3
$rng = $ws.Range("C3,E3,G3,I3,K3,M3,O3")
try {
  $ser.ErrorBar(2, 1, 2, $rng)
  Write-Host "ErrorBar with range worked"
} catch {
  Write-Host "ErrorBar with range failed:" $_.Exception.Message
}
